$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as text (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.936.23"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "1.912.81"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "325.20"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4593"
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").Value = "0.07732"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "0.9807"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "22.10"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("D12").Value = "1.921.32"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "6.949"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "5.672"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "0.07045"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "83.87"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "0.000009473"
$ws.Range("E18").Value = "  -3.45%  "
$ws.Range("D19").Value = "16.69"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "28.910.45"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("D23").Value = "10.89"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").Value = "2.092"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "158.56"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "19.01"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "5.668"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "117.40"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "1.861"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.09293"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "0.8700"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "5.084"
$ws.Range("E32").Value = "  -1.87%  "
$ws.Range("D33").Value = "1.253"
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "3.151"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "0.05729"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "1.165"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("D38").Value = "0.02046"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "7.413"
$ws.Range("E39").Value = "  -3.49%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.5492"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").Value = "0.1756"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "2.860"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("D43").Value = "9.321"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "0.5185"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "11.26"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "0.06916"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "2.101"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "0.000002579"
$ws.Range("E48").Value = "  -8.17%  "
$ws.Range("D49").Value = "1.781"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").Value = "110.55"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "0.2879"
$ws.Range("E51").Value = "  -4.31%  "
